$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-46 (B and C columns)
$ws.Cells.Item(2, 2).Value = 2.163010895737672
$ws.Cells.Item(2, 3).Value = 4.097185091623455
$ws.Cells.Item(3, 2).Value = 3.235570498670756
$ws.Cells.Item(3, 3).Value = 8.192073022743907
$ws.Cells.Item(4, 2).Value = 4.765536510082268
$ws.Cells.Item(4, 3).Value = 12.40482765156577
$ws.Cells.Item(5, 2).Value = 6.433122933805482
$ws.Cells.Item(5, 3).Value = 17.00568812445333
$ws.Cells.Item(6, 2).Value = 8.244538086024699
$ws.Cells.Item(6, 3).Value = 21.20644016298202
$ws.Cells.Item(7, 2).Value = 9.716845088336667
$ws.Cells.Item(7, 3).Value = 25.41009155653449
$ws.Cells.Item(8, 2).Value = 15.10078209939159
$ws.Cells.Item(8, 3).Value = 29.48681143145043
$ws.Cells.Item(9, 2).Value = 15.60091426436921
$ws.Cells.Item(9, 3).Value = 33.94457234743799
$ws.Cells.Item(10, 2).Value = 18.28937490691169
$ws.Cells.Item(10, 3).Value = 38.18049870304964
$ws.Cells.Item(11, 2).Value = 19.78273921839277
$ws.Cells.Item(11, 3).Value = 42.85543957863415
$ws.Cells.Item(12, 2).Value = 20.49510616167762
$ws.Cells.Item(12, 3).Value = 47.22440470670529
$ws.Cells.Item(13, 2).Value = 22.35895415147025
$ws.Cells.Item(13, 3).Value = 51.39020283779726
$ws.Cells.Item(14, 2).Value = 23.02380353402351
$ws.Cells.Item(14, 3).Value = 55.68967679406274
$ws.Cells.Item(15, 2).Value = 23.87020025854967
$ws.Cells.Item(15, 3).Value = 60.53698403023608
$ws.Cells.Item(16, 2).Value = 25.62639097014731
$ws.Cells.Item(16, 3).Value = 65.11865146818074
$ws.Cells.Item(17, 2).Value = 27.68737114814559
$ws.Cells.Item(17, 3).Value = 69.35398882958546
$ws.Cells.Item(18, 2).Value = 28.04081588552889
$ws.Cells.Item(18, 3).Value = 73.86874721055166
$ws.Cells.Item(19, 2).Value = 30.70183044595878
$ws.Cells.Item(19, 3).Value = 78.20554293012546
$ws.Cells.Item(20, 2).Value = 31.60273437940791
$ws.Cells.Item(20, 3).Value = 82.3801181274232
$ws.Cells.Item(21, 2).Value = 35.53007301666943
$ws.Cells.Item(21, 3).Value = 87.17996028810926
$ws.Cells.Item(22, 2).Value = 37.6003769054215
$ws.Cells.Item(22, 3).Value = 91.52613326647344
$ws.Cells.Item(23, 2).Value = 39.05853914299175
$ws.Cells.Item(23, 3).Value = 96.84302110174271
$ws.Cells.Item(24, 2).Value = 40.76999755742614
$ws.Cells.Item(24, 3).Value = 101.3558014702302
$ws.Cells.Item(25, 2).Value = 41.70417525617449
$ws.Cells.Item(25, 3).Value = 105.5792061782314
$ws.Cells.Item(26, 2).Value = 45.38490428784193
$ws.Cells.Item(26, 3).Value = 109.7988444645075
$ws.Cells.Item(27, 2).Value = 47.82535146328136
$ws.Cells.Item(27, 3).Value = 114.1942301063541
$ws.Cells.Item(28, 2).Value = 48.98595736844501
$ws.Cells.Item(28, 3).Value = 118.7456038863531
$ws.Cells.Item(29, 2).Value = 52.37104685418971
$ws.Cells.Item(29, 3).Value = 123.143117136389
$ws.Cells.Item(30, 2).Value = 53.18489340030789
$ws.Cells.Item(30, 3).Value = 127.9731537246993
$ws.Cells.Item(31, 2).Value = 54.1711802627469
$ws.Cells.Item(31, 3).Value = 132.2565220342947
$ws.Cells.Item(32, 2).Value = 54.73512696995507
$ws.Cells.Item(32, 3).Value = 136.5374431658016
$ws.Cells.Item(33, 2).Value = 56.38404311845105
$ws.Cells.Item(33, 3).Value = 141.3111997568723
$ws.Cells.Item(34, 2).Value = 59.41643852151996
$ws.Cells.Item(34, 3).Value = 145.4429860812598
$ws.Cells.Item(35, 2).Value = 62.17543693949342
$ws.Cells.Item(35, 3).Value = 150.4380438311324
$ws.Cells.Item(36, 2).Value = 65.22913225605791
$ws.Cells.Item(36, 3).Value = 154.7707987036789
$ws.Cells.Item(37, 2).Value = 66.73505044191529
$ws.Cells.Item(37, 3).Value = 159.1314179269283
$ws.Cells.Item(38, 2).Value = 70.13988725944452
$ws.Cells.Item(38, 3).Value = 163.309299356001
$ws.Cells.Item(39, 2).Value = 71.45184723680686
$ws.Cells.Item(39, 3).Value = 167.7036469500569
$ws.Cells.Item(40, 2).Value = 72.77775168105482
$ws.Cells.Item(40, 3).Value = 171.8206072148614
$ws.Cells.Item(41, 2).Value = 75.70011118403229
$ws.Cells.Item(41, 3).Value = 176.0209181640597
$ws.Cells.Item(42, 2).Value = 78.04472906015513
$ws.Cells.Item(42, 3).Value = 180.3027379676512
$ws.Cells.Item(43, 2).Value = 80.60745366192096
$ws.Cells.Item(43, 3).Value = 184.8225738964389
$ws.Cells.Item(44, 2).Value = 82.41558936556531
$ws.Cells.Item(44, 3).Value = 189.3561949000019
$ws.Cells.Item(45, 2).Value = 84.36098443395464
$ws.Cells.Item(45, 3).Value = 194.0112015046877
$ws.Cells.Item(46, 2).Value = 85.94822765197023
$ws.Cells.Item(46, 3).Value = 198.1989190527542

# Add new rows 47-50
$ws.Cells.Item(47, 1).Value = 45
$ws.Cells.Item(46, 1).Copy() | Out-Null
$ws.Cells.Item(47, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(47, 2).Value = 90.43212939462663
$ws.Cells.Item(47, 3).Value = 202.7212422050105
$ws.Cells.Item(48, 1).Value = 46
$ws.Cells.Item(46, 1).Copy() | Out-Null
$ws.Cells.Item(48, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(48, 2).Value = 94.47960076638105
$ws.Cells.Item(48, 3).Value = 207.4422852410104
$ws.Cells.Item(49, 1).Value = 47
$ws.Cells.Item(46, 1).Copy() | Out-Null
$ws.Cells.Item(49, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(49, 2).Value = 97.75776272255001
$ws.Cells.Item(49, 3).Value = 211.645803408932
$ws.Cells.Item(50, 1).Value = 48
$ws.Cells.Item(46, 1).Copy() | Out-Null
$ws.Cells.Item(50, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(50, 2).Value = 99.86844980236503
$ws.Cells.Item(50, 3).Value = 215.8086895891753
